# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos - refreshed case counts ---
$ws.Range("B4").Value = 710021
$ws.Range("C4").Value = 286
$ws.Range("E4").Value = 612353
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 37158

# --- Row 10: China - refreshed case counts ---
$ws.Range("C10").Value = 27
$ws.Range("G10").Value = 0

# --- Rows 37/38: Chequia & Australia swap places (Australia's totals
#     overtook Chequia's), and Australia gets fresh data while
#     Chequia's figures stay the same but move down one row ---
$ws.Range("A37").Value = "Australia"
$ws.Range("B37").Value = 6560
$ws.Range("C37").Value = 27
$ws.Range("D37").Value = 3821
$ws.Range("E37").Value = 2672
$ws.Range("F37").Value = 57
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 67

$ws.Range("A38").Value = "Chequia"
$ws.Range("B38").Value = 6549
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 1174
$ws.Range("E38").Value = 5202
$ws.Range("F38").Value = 82
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 173
